$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The sheet previously held 12 rows of "undirected arc" pairs, each cell
# formatted with the built-in "Moderate" (yellow/orange) cell style. The
# new data is just 3 rows, written as plain (unformatted) numbers, so drop
# the old content/format first.
$ws.Range("A1:B12").Clear()

# That custom cell style is no longer used anywhere in the workbook, so
# remove it from the workbook's style gallery too.
$wb.Styles.Item(2).Delete()

# New set of undirected arcs.
$data = @(
    @(2, 4),
    @(10, 11),
    @(30, 31)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Match the saved view state: columns A:B selected, active cell in B1.
$ws.Range("A1:B1048576").Select()
